# إضافة حدث جديد في Card12 by admin at 2026-02-18 13:19:08
#
# Row 32 previously had blank placeholder cells (B32:K32, P32) that now get
# the literal text "nan" (matching the rest of the sheet's "no data" marker),
# and a brand-new row 33 is appended with the newest service-log entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# --- Row 32: fill the previously-empty columns with the "nan" placeholder ---
$ws.Range("B32:K32").Value = "nan"
$ws.Range("P32").Value = "nan"

# --- Row 33: new service event row ---
# Column A holds the card number as text ("12"), same as every other row in
# this sheet, so force text using a quote-prefix to avoid Excel's automatic
# number coercion.
$ws.Range("A33").Value = "'12"

# Columns B:K and P have no recorded value for this event - keep them present
# as empty text cells (same shape as the rest of the sheet), not plain blanks.
$ws.Range("B33:K33").Value = "'"
$ws.Range("P33").Value = "'"

$ws.Range("L33").Value = "17/2/2026"
$ws.Range("M33").Value = "حسام و ايهاب"
$ws.Range("N33").Value = "تجربه"
$ws.Range("O33").Value = "تم تغيير  سير 1200 الدوفر المسنن للتجربه"
